# data-grid-sample.xlsx edit
# - rename the sheet
# - replace the 2-row sample data with a 4-row grid
# - narrow the first three columns
# - move the selection back to A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "data-grid-sample"

# Row 1 (header-ish row)
$ws.Range("A1").Value = "first"
$ws.Range("B1").Value = "things"
$ws.Range("C1").Value = "phirst"

# Row 2
$ws.Range("A2").Value = '"derp"'
$ws.Range("B2").Value = "is"
$ws.Range("C2").Value = "a,word"

# Row 3 (the old formula cell becomes a plain literal number)
$ws.Range("A3").Value = "over"
$ws.Range("B3").Value = 9000
$ws.Range("C3").Value = "duhh"

# Row 4 (new)
$ws.Range("A4").Value = "magic"
$ws.Range("B4").Value = "fourth"
$ws.Range("C4").Value = "row"

# Narrow the first three columns
$ws.Columns.Item(1).ColumnWidth = 5.417
$ws.Columns.Item(2).ColumnWidth = 5.25
$ws.Columns.Item(3).ColumnWidth = 4.917

# Reset the selection to A1
$null = $ws.Range("A1").Select()
